$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (Price -> column D,
# Volume(1h) -> column E) and apply the Fetch.AI/Monero row swap (rows 37-38)
# described by the upstream diff. Price cells are plain text in the workbook
# (e.g. thousand-separated with dots, or fixed-decimal strings like '143.80'),
# so numeric-looking values are written with a leading apostrophe to keep them
# stored as text instead of being auto-converted to numbers by Excel.

$ws.Range('D2').Value = '65.731.14'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '3.505.24'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'596.23"
$ws.Range('E5').Value = '  -1.50%  '
$ws.Range('D6').Value = "'143.80"
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('D7').Value = '3.499.83'
$ws.Range('E7').Value = '  -1.43%  '
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('D9').Value = "'0.496"
$ws.Range('E9').Value = '  +0.55%  '
$ws.Range('D10').Value = "'0.133"
$ws.Range('E10').Value = '  -2.13%  '
$ws.Range('D11').Value = "'7.65"
$ws.Range('E11').Value = '  -3.26%  '
$ws.Range('D12').Value = "'0.403"
$ws.Range('E12').Value = '  -2.48%  '
$ws.Range('D13').Value = '4.105.01'
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('D14').Value = "'0.0000199"
$ws.Range('E14').Value = '  -4.01%  '
$ws.Range('D15').Value = "'28.62"
$ws.Range('E15').Value = '  -4.62%  '
$ws.Range('D16').Value = '3.511.31'
$ws.Range('E16').Value = '  -1.29%  '
$ws.Range('E17').Value = '  +1.17%  '
$ws.Range('D18').Value = '65.791.06'
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('D19').Value = "'10.93"
$ws.Range('E19').Value = '  -5.21%  '
$ws.Range('D20').Value = "'6.16"
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('D21').Value = "'14.32"
$ws.Range('E21').Value = '  -3.55%  '
$ws.Range('D22').Value = "'412.88"
$ws.Range('E22').Value = '  -4.27%  '
$ws.Range('D23').Value = "'0.594"
$ws.Range('E23').Value = '  -2.60%  '
$ws.Range('D24').Value = "'77.48"
$ws.Range('E24').Value = '  -2.76%  '
$ws.Range('D25').Value = '3.655.26'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -4.14%  '
$ws.Range('D28').Value = "'8.99"
$ws.Range('E28').Value = '  -2.25%  '
$ws.Range('D29').Value = "'2.43"
$ws.Range('E29').Value = '  -2.99%  '
$ws.Range('D30').Value = "'7.65"
$ws.Range('E30').Value = '  -4.01%  '
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('D32').Value = '3.507.72'
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('D33').Value = "'0.152"
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('D34').Value = "'24.13"
$ws.Range('E34').Value = '  -4.94%  '
$ws.Range('D36').Value = "'7.45"
$ws.Range('E36').Value = '  -5.36%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = "'174.15"
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').Value = "'1.24"
$ws.Range('E38').Value = '  -14.75%  '
$ws.Range('D39').Value = "'5.19"
$ws.Range('E39').Value = '  -7.27%  '
$ws.Range('D40').Value = "'1.56"
$ws.Range('E40').Value = '  -9.34%  '
$ws.Range('D41').Value = "'0.0815"
$ws.Range('E41').Value = '  -3.87%  '
$ws.Range('D42').Value = "'5.02"
$ws.Range('E42').Value = '  -3.33%  '
$ws.Range('D43').Value = "'0.850"
$ws.Range('E43').Value = '  -4.24%  '
$ws.Range('D44').Value = "'45.15"
$ws.Range('E44').Value = '  -1.98%  '
$ws.Range('D45').Value = "'1.76"
$ws.Range('E45').Value = '  -8.50%  '
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('D47').Value = "'2.39"
$ws.Range('E47').Value = '  -5.38%  '
$ws.Range('D48').Value = "'7.04"
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('D49').Value = "'22.40"
$ws.Range('E49').Value = '  -4.48%  '
$ws.Range('E50').Value = '  -8.93%  '
$ws.Range('D51').Value = "'22.72"
$ws.Range('E51').Value = '  -9.28%  '
